$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4640.4614
$ws.Range("I40").Value = 6277.8887
$ws.Range("J40").Value = 956.25
$ws.Range("K40").Value = 6277.8887
$ws.Range("L40").Value = 956.25
$ws.Range("M40").Value = -6102.8887
$ws.Range("N40").Value = -1306.25

$ws.Range("H64").Value = 23031.96
$ws.Range("I64").Value = 3104.95
$ws.Range("J64").Value = 102740
$ws.Range("K64").Value = 3104.95
$ws.Range("L64").Value = 102740
$ws.Range("M64").Value = -2856.95
$ws.Range("N64").Value = -103236

$ws.Range("H67").Value = 23031.96
$ws.Range("I67").Value = 3104.95
$ws.Range("J67").Value = 102740
$ws.Range("K67").Value = 3104.95
$ws.Range("L67").Value = 102740
$ws.Range("M67").Value = -2246.95
$ws.Range("N67").Value = -104456

$ws.Range("H74").Value = 2130695.5
$ws.Range("I74").Value = 2383961.2
$ws.Range("K74").Value = 2383961.2
$ws.Range("M74").Value = -2383025.2

$ws.Range("H77").Value = 2130695.5
$ws.Range("I77").Value = 2383961.2
$ws.Range("K77").Value = 11919806
$ws.Range("M77").Value = -11915126

$ws.Range("H129").Value = 608.8570999999999
$ws.Range("I129").Value = 502.18182
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 1506.54546
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 3493.45454
$ws.Range("N129").Value = -13000

$ws.Range("H132").Value = 138126.94
$ws.Range("I132").Value = 3087
$ws.Range("K132").Value = 9261
$ws.Range("M132").Value = -6731

$ws.Range("H137").Value = 36533.195
$ws.Range("J137").Value = 6605.4443
$ws.Range("L137").Value = 19816.3329
$ws.Range("N137").Value = -24916.3329

$ws.Range("H138").Value = 1717.129
$ws.Range("I138").Value = 827.9756
$ws.Range("J138").Value = 2418.1924
$ws.Range("K138").Value = 2483.9268
$ws.Range("L138").Value = 7254.5772
$ws.Range("M138").Value = 2656.0732
$ws.Range("N138").Value = -17534.5772

$ws.Range("H141").Value = 1674.3125
$ws.Range("I141").Value = 1068.9
$ws.Range("J141").Value = 2683.3333
$ws.Range("K141").Value = 3206.7
$ws.Range("L141").Value = 8049.999899999999
$ws.Range("M141").Value = 1973.3
$ws.Range("N141").Value = -18409.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 259.18182
$ws.Range("I5").Value = 185.1
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 185.1
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = -73.09999999999999
$ws.Range("N5").Value = -1224

$ws.Range("H37").Value = 9263.916999999999
$ws.Range("I37").Value = 3056
$ws.Range("J37").Value = 11333.223
$ws.Range("K37").Value = 3056
$ws.Range("L37").Value = 11333.223
$ws.Range("M37").Value = -2783
$ws.Range("N37").Value = -11879.223

$ws.Range("H44").Value = 22227.666
$ws.Range("J44").Value = 22227.666
$ws.Range("L44").Value = 22227.666
$ws.Range("N44").Value = -23203.666

$ws.Range("H49").Value = 5000
$ws.Range("J49").Value = 5000
$ws.Range("L49").Value = 5000
$ws.Range("N49").Value = -5520

$ws.Range("H55").Value = 13285.714
$ws.Range("J55").Value = 13285.714
$ws.Range("L55").Value = 13285.714
$ws.Range("N55").Value = -13915.714

$ws.Range("H59").Value = 20000
$ws.Range("I59").Value = 20000
$ws.Range("K59").Value = 20000
$ws.Range("M59").Value = -19196

$ws.Range("H63").Value = 2998.6
$ws.Range("I63").Value = 2998.6
$ws.Range("K63").Value = 2998.6
$ws.Range("M63").Value = -2312.6

$ws.Range("H66").Value = 2998.6
$ws.Range("I66").Value = 2998.6
$ws.Range("K66").Value = 14993
$ws.Range("M66").Value = -11561

$ws.Range("H80").Value = 16358.5
$ws.Range("I80").Value = 5001
$ws.Range("J80").Value = 17232.154
$ws.Range("K80").Value = 5001
$ws.Range("L80").Value = 17232.154
$ws.Range("M80").Value = -4003
$ws.Range("N80").Value = -19228.154

$ws.Range("H83").Value = 16358.5
$ws.Range("I83").Value = 5001
$ws.Range("J83").Value = 17232.154
$ws.Range("K83").Value = 15003
$ws.Range("L83").Value = 51696.462
$ws.Range("M83").Value = -10011
$ws.Range("N83").Value = -61680.462

$ws.Range("H88").Value = 1800
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 1800
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = $null

$ws.Range("H91").Value = 1800
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 1800
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = $null

$ws.Range("H102").Value = 1966.6666
$ws.Range("I102").Value = 1860
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1860
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -238
$ws.Range("N102").Value = -5744

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 259.18182
$ws.Range("I4").Value = 185.1
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 185.1
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -70.09999999999999
$ws.Range("N4").Value = -1230

$ws.Range("H15").Value = 18695.8
$ws.Range("J15").Value = 18695.8
$ws.Range("L15").Value = 18695.8
$ws.Range("N15").Value = -19149.8

$ws.Range("H35").Value = 17037
$ws.Range("J35").Value = 17037
$ws.Range("L35").Value = 17037
$ws.Range("N35").Value = -17657

$ws.Range("H86").Value = 319963.72
$ws.Range("I86").Value = 1824.5834
$ws.Range("J86").Value = 701730.7
$ws.Range("K86").Value = 1824.5834
$ws.Range("L86").Value = 701730.7
$ws.Range("M86").Value = -701.5834
$ws.Range("N86").Value = -703976.7

$ws.Range("H89").Value = 319963.72
$ws.Range("I89").Value = 1824.5834
$ws.Range("J89").Value = 701730.7
$ws.Range("K89").Value = 9122.916999999999
$ws.Range("L89").Value = 3508653.5
$ws.Range("M89").Value = -3506.916999999999
$ws.Range("N89").Value = -3519885.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1066.875
$ws.Range("I16").Value = 1007.3333
$ws.Range("J16").Value = 1102.6
$ws.Range("K16").Value = 1007.3333
$ws.Range("L16").Value = 1102.6
$ws.Range("M16").Value = -720.3333
$ws.Range("N16").Value = -1676.6

$ws.Range("H22").Value = 1149
$ws.Range("I22").Value = 1415.3334
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 1415.3334
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -1065.3334
$ws.Range("N22").Value = -1050

$ws.Range("H31").Value = 9094.654
$ws.Range("I31").Value = 7420.579
$ws.Range("J31").Value = 13638.571
$ws.Range("K31").Value = 7420.579
$ws.Range("L31").Value = 13638.571
$ws.Range("M31").Value = -7125.579
$ws.Range("N31").Value = -14228.571

$ws.Range("H34").Value = 9094.654
$ws.Range("I34").Value = 7420.579
$ws.Range("J34").Value = 13638.571
$ws.Range("K34").Value = 7420.579
$ws.Range("L34").Value = 13638.571
$ws.Range("M34").Value = -7218.579
$ws.Range("N34").Value = -14042.571

$ws.Range("H113").Value = 1066.875
$ws.Range("I113").Value = 1007.3333
$ws.Range("J113").Value = 1102.6
$ws.Range("K113").Value = 1007.3333
$ws.Range("L113").Value = 1102.6
$ws.Range("M113").Value = 1162.6667
$ws.Range("N113").Value = -5442.6

$ws.Range("H132").Value = 818.119
$ws.Range("I132").Value = 719.4722
$ws.Range("J132").Value = 1410
$ws.Range("K132").Value = 2158.4166
$ws.Range("L132").Value = 4230
$ws.Range("M132").Value = 371.5834
$ws.Range("N132").Value = -9290

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2003
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 2003
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 6009
$ws.Range("M49").Value = $null
$ws.Range("N49").Value = -6321

$ws.Range("H122").Value = 375.25714
$ws.Range("I122").Value = 312.8421
$ws.Range("J122").Value = 449.375
$ws.Range("K122").Value = 2815.5789
$ws.Range("L122").Value = 4044.375
$ws.Range("M122").Value = -365.5789
$ws.Range("N122").Value = -8944.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 13699.538
$ws.Range("I9").Value = 182.85715
$ws.Range("J9").Value = 29469
$ws.Range("K9").Value = 182.85715
$ws.Range("L9").Value = 29469
$ws.Range("M9").Value = 41.14285000000001
$ws.Range("N9").Value = -29917

$ws.Range("H132").Value = 506042.12
$ws.Range("I132").Value = 205298.7
$ws.Range("J132").Value = 720858.9
$ws.Range("K132").Value = 615896.1000000001
$ws.Range("L132").Value = 2162576.7
$ws.Range("M132").Value = -613366.1000000001
$ws.Range("N132").Value = -2167636.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 863.7143
$ws.Range("I126").Value = 736.6
$ws.Range("J126").Value = 1181.5
$ws.Range("K126").Value = 2209.8
$ws.Range("L126").Value = 3544.5
$ws.Range("M126").Value = 260.1999999999998
$ws.Range("N126").Value = -8484.5
